$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "43.883.77"
$ws.Range("E2").Value = "  -0.66%  "

# Row 3
$ws.Range("D3").Value = "2.346.50"
$ws.Range("E3").Value = "  -2.17%  "

# Row 4
$ws.Range("E4").Value = "  -0.02%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.81"
$ws.Range("E5").Value = "  -1.90%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.669"
$ws.Range("E6").Value = "  -4.52%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "72.37"
$ws.Range("E7").Value = "  -6.40%  "

# Row 8
$ws.Range("E8").Value = "  -0.03%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.589"
$ws.Range("E9").Value = "  -9.00%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0996"
$ws.Range("E10").Value = "  -4.07%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "58.48"
$ws.Range("E11").Value = "  +1.02%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "32.44"
$ws.Range("E12").Value = "  -4.31%  "

# Row 13
$ws.Range("E13").Value = "  -0.77%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.16"
$ws.Range("E14").Value = "  -6.54%  "

# Row 15
$ws.Range("D15").Value = "2.695.92"
$ws.Range("E15").Value = "  -2.30%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "16.25"
$ws.Range("E16").Value = "  -5.70%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.901"
$ws.Range("E17").Value = "  -4.66%  "

# Row 18
$ws.Range("D18").Value = "2.354.37"
$ws.Range("E18").Value = "  -2.11%  "

# Row 19
$ws.Range("D19").Value = "43.784.05"
$ws.Range("E19").Value = "  -0.95%  "

# Row 20
$ws.Range("E20").Value = "  -1.92%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.66"
$ws.Range("E21").Value = "  -1.50%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "77.81"
$ws.Range("E22").Value = "  -1.89%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "253.76"
$ws.Range("E23").Value = "  -2.52%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.92"
$ws.Range("E24").Value = "  +7.12%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("E25").Value = "  +0.00%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.72"
$ws.Range("E26").Value = "  +1.08%  "

# Row 27
$ws.Range("E27").Value = "  -2.78%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.39"
$ws.Range("E28").Value = "  -7.16%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.26"
$ws.Range("E29").Value = "  -2.20%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "177.08"
$ws.Range("E30").Value = "  +0.61%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "22.29"
$ws.Range("E31").Value = "  -5.57%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.127"
$ws.Range("E32").Value = "  -3.43%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.135"
$ws.Range("E33").Value = "  -1.16%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0740"
$ws.Range("E34").Value = "  -3.30%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.08"
$ws.Range("E35").Value = "  -6.74%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.37"
$ws.Range("E36").Value = "  -1.26%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.73"
$ws.Range("E37").Value = "  -4.83%  "

# Row 38
$ws.Range("B38").Value = "THORChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.41"
$ws.Range("E38").Value = "  -3.43%  "

# Row 39
$ws.Range("B39").Value = "LidoDAOToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.38"
$ws.Range("E39").Value = "  -5.71%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0273"
$ws.Range("E40").Value = "  -1.99%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "64.53"
$ws.Range("E41").Value = "  +15.93%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.14"
$ws.Range("E42").Value = "  +12.33%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "9.21"
$ws.Range("E43").Value = "  +1.01%  "

# Row 44
$ws.Range("E44").Value = "  +5.42%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "18.72"
$ws.Range("E45").Value = "  -2.16%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.198"
$ws.Range("E46").Value = "  -3.27%  "

# Row 47
$ws.Range("E47").Value = "  -0.02%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.46"
$ws.Range("E48").Value = "  -4.89%  "

# Row 49
$ws.Range("E49").Value = "  -4.50%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "98.66"
$ws.Range("E50").Value = "  -6.04%  "

# Row 51
$ws.Range("E51").Value = "  -7.22%  "
